$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2, J2, AG2, AH2, AI2)
$ws.Range("A2").Value = 2794
$ws.Range("J2").Value = 365
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 212239.4
$ws.Range("AI2").Value = 222760.6

# Row 3 (A3, J3)
$ws.Range("A3").Value = 2793
$ws.Range("J3").Value = 365

# Row 4 (A4, J4)
$ws.Range("A4").Value = 3010
$ws.Range("J4").Value = 351

# Row 5 (A5, J5)
$ws.Range("A5").Value = 3017
$ws.Range("J5").Value = 379

# Row 6 (A6, J6)
$ws.Range("A6").Value = 2993
$ws.Range("J6").Value = 434

# Row 7 (A7, J7)
$ws.Range("A7").Value = 3287
$ws.Range("J7").Value = 379

# Row 8 (A8, J8, AG8, AH8, AI8)
$ws.Range("A8").Value = 3028
$ws.Range("J8").Value = 434
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0.04
$ws.Range("AI8").Value = 20799.96

# Row 9 (A9, J9, AG9, AH9, AI9)
$ws.Range("A9").Value = 3027
$ws.Range("J9").Value = 434
$ws.Range("AG9").Value = 0
$ws.Range("AH9").Value = 0.02
$ws.Range("AI9").Value = 33259.980000000003
